$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 4.861952666666666
$ws.Cells.Item(2, 8).Value = 14.585858
$ws.Cells.Item(2, 9).Value = 0.3995648519435639
$ws.Cells.Item(2, 10).Value = 0.3995648519435638
$ws.Cells.Item(2, 15).Value = 0.9347132976570145
$ws.Cells.Item(2, 16).Value = 0.9347132976570145
$ws.Cells.Item(2, 17).Value = 43.18741443143089
$ws.Cells.Item(2, 18).Value = 388.6867298828779
$ws.Cells.Item(2, 19).Value = 0.3734785803880054
$ws.Cells.Item(2, 20).Value = 0.3734785803880053

$ws.Cells.Item(3, 7).Value = 4.861952666666666
$ws.Cells.Item(3, 8).Value = 14.585858
$ws.Cells.Item(3, 9).Value = 0.3995648519435639
$ws.Cells.Item(3, 10).Value = 0.3995648519435638
$ws.Cells.Item(3, 13).Value = 0.616144
$ws.Cells.Item(3, 14).Value = 1.848432
$ws.Cells.Item(3, 15).Value = 0.06483569448352988
$ws.Cells.Item(3, 16).Value = 0.0648356944835299
$ws.Cells.Item(3, 17).Value = 2.995662963850667
$ws.Cells.Item(3, 18).Value = 26.960966674656
$ws.Cells.Item(3, 19).Value = 0.02590606466696976
$ws.Cells.Item(3, 20).Value = 0.02590606466696976

$ws.Cells.Item(4, 7).Value = 4.861952666666666
$ws.Cells.Item(4, 8).Value = 14.585858
$ws.Cells.Item(4, 9).Value = 0.3995648519435639
$ws.Cells.Item(4, 10).Value = 0.3995648519435638
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.004286
$ws.Cells.Item(4, 14).Value = 0.012858
$ws.Cells.Item(4, 15).Value = 0.0004510078594555965
$ws.Cells.Item(4, 16).Value = 0.0004510078594555965
$ws.Cells.Item(4, 17).Value = 0.02083832912933333
$ws.Cells.Item(4, 18).Value = 0.187544962164
$ws.Cells.Item(4, 19).Value = 0.0001802068885887591
$ws.Cells.Item(4, 20).Value = 0.0001802068885887591

$ws.Cells.Item(5, 9).Value = 0.04932556406896855
$ws.Cells.Item(5, 10).Value = 0.04932556406896854
$ws.Cells.Item(5, 15).Value = 0.9347132976570145
$ws.Cells.Item(5, 16).Value = 0.9347132976570145
$ws.Cells.Item(5, 18).Value = 47.982679418218
$ws.Cells.Item(5, 19).Value = 0.04610526064969794
$ws.Cells.Item(5, 20).Value = 0.04610526064969793

$ws.Cells.Item(6, 9).Value = 0.04932556406896855
$ws.Cells.Item(6, 10).Value = 0.04932556406896854
$ws.Cells.Item(6, 13).Value = 0.616144
$ws.Cells.Item(6, 14).Value = 1.848432
$ws.Cells.Item(6, 15).Value = 0.06483569448352988
$ws.Cells.Item(6, 16).Value = 0.0648356944835299
$ws.Cells.Item(6, 17).Value = 0.3698092180373334
$ws.Cells.Item(6, 18).Value = 3.328282962336
$ws.Cells.Item(6, 19).Value = 0.003198057202203424
$ws.Cells.Item(6, 20).Value = 0.003198057202203424

$ws.Cells.Item(7, 9).Value = 0.04932556406896855
$ws.Cells.Item(7, 10).Value = 0.04932556406896854
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.004286
$ws.Cells.Item(7, 14).Value = 0.012858
$ws.Cells.Item(7, 15).Value = 0.0004510078594555965
$ws.Cells.Item(7, 16).Value = 0.0004510078594555965
$ws.Cells.Item(7, 17).Value = 0.002572454342666667
$ws.Cells.Item(7, 18).Value = 0.023152089084
$ws.Cells.Item(7, 19).Value = 0.00002224621706718539
$ws.Cells.Item(7, 20).Value = 0.00002224621706718539

$ws.Cells.Item(8, 7).Value = 4.206754333333333
$ws.Cells.Item(8, 8).Value = 12.620263
$ws.Cells.Item(8, 9).Value = 0.3457193616641432
$ws.Cells.Item(8, 10).Value = 0.3457193616641432
$ws.Cells.Item(8, 15).Value = 0.9347132976570145
$ws.Cells.Item(8, 16).Value = 0.9347132976570145
$ws.Cells.Item(8, 17).Value = 37.36746432158144
$ws.Cells.Item(8, 18).Value = 336.3071788942329
$ws.Cells.Item(8, 19).Value = 0.3231484846049694
$ws.Cells.Item(8, 20).Value = 0.3231484846049694

$ws.Cells.Item(9, 7).Value = 4.206754333333333
$ws.Cells.Item(9, 8).Value = 12.620263
$ws.Cells.Item(9, 9).Value = 0.3457193616641432
$ws.Cells.Item(9, 10).Value = 0.3457193616641432
$ws.Cells.Item(9, 13).Value = 0.616144
$ws.Cells.Item(9, 14).Value = 1.848432
$ws.Cells.Item(9, 15).Value = 0.06483569448352988
$ws.Cells.Item(9, 16).Value = 0.0648356944835299
$ws.Cells.Item(9, 17).Value = 2.591966441957333
$ws.Cells.Item(9, 18).Value = 23.327697977616
$ws.Cells.Item(9, 19).Value = 0.02241495490989736
$ws.Cells.Item(9, 20).Value = 0.02241495490989737

$ws.Cells.Item(10, 7).Value = 4.206754333333333
$ws.Cells.Item(10, 8).Value = 12.620263
$ws.Cells.Item(10, 9).Value = 0.3457193616641432
$ws.Cells.Item(10, 10).Value = 0.3457193616641432
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.004286
$ws.Cells.Item(10, 14).Value = 0.012858
$ws.Cells.Item(10, 15).Value = 0.0004510078594555965
$ws.Cells.Item(10, 16).Value = 0.0004510078594555965
$ws.Cells.Item(10, 17).Value = 0.01803014907266666
$ws.Cells.Item(10, 18).Value = 0.162271341654
$ws.Cells.Item(10, 19).Value = 0.0001559221492765004
$ws.Cells.Item(10, 20).Value = 0.0001559221492765005

$ws.Cells.Item(11, 7).Value = 2.499212666666667
$ws.Cells.Item(11, 8).Value = 7.497638
$ws.Cells.Item(11, 9).Value = 0.2053902223233243
$ws.Cells.Item(11, 10).Value = 0.2053902223233243
$ws.Cells.Item(11, 15).Value = 0.9347132976570145
$ws.Cells.Item(11, 16).Value = 0.9347132976570145
$ws.Cells.Item(11, 17).Value = 22.19983216365089
$ws.Cells.Item(11, 18).Value = 199.798489472858
$ws.Cells.Item(11, 19).Value = 0.1919809720143418
$ws.Cells.Item(11, 20).Value = 0.1919809720143418

$ws.Cells.Item(12, 7).Value = 2.499212666666667
$ws.Cells.Item(12, 8).Value = 7.497638
$ws.Cells.Item(12, 9).Value = 0.2053902223233243
$ws.Cells.Item(12, 10).Value = 0.2053902223233243
$ws.Cells.Item(12, 13).Value = 0.616144
$ws.Cells.Item(12, 14).Value = 1.848432
$ws.Cells.Item(12, 15).Value = 0.06483569448352988
$ws.Cells.Item(12, 16).Value = 0.0648356944835299
$ws.Cells.Item(12, 17).Value = 1.539874889290667
$ws.Cells.Item(12, 18).Value = 13.858874003616
$ws.Cells.Item(12, 19).Value = 0.01331661770445934
$ws.Cells.Item(12, 20).Value = 0.01331661770445934

$ws.Cells.Item(13, 7).Value = 2.499212666666667
$ws.Cells.Item(13, 8).Value = 7.497638
$ws.Cells.Item(13, 9).Value = 0.2053902223233243
$ws.Cells.Item(13, 10).Value = 0.2053902223233243
$ws.Cells.Item(13, 11).Value = 2
$ws.Cells.Item(13, 12).Value = 0.6666666666666666
$ws.Cells.Item(13, 13).Value = 0.004286
$ws.Cells.Item(13, 14).Value = 0.012858
$ws.Cells.Item(13, 15).Value = 0.0004510078594555965
$ws.Cells.Item(13, 16).Value = 0.0004510078594555965
$ws.Cells.Item(13, 17).Value = 0.01071162548933333
$ws.Cells.Item(13, 18).Value = 0.09640462940400001
$ws.Cells.Item(13, 19).Value = 0.00009263260452315158
$ws.Cells.Item(13, 20).Value = 0.00009263260452315158
